{"js": "// The document currently starts with a paragraph \"OBSTACLES: Lerp horizontal\n// and vertical, lerp scale too\" followed by an empty paragraph, etc.\n// The edit:\n//   1) Inserts a brand-new paragraph \"ADD GAME MUSIC\" before everything else.\n//   2) Inserts a brand-new paragraph (with a mix of plain/bold runs) right\n//      after it, before the old first paragraph.\n//   3) Leaves the rest of the body alone, except the final\n//      \"Recognized Voice Commands\" paragraph gets extra text spliced in.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst originalFirstPara = body.paragraphs.items[0];\n\n// 1) New paragraph: \"ADD GAME MUSIC\"\nconst titlePara = originalFirstPara.insertParagraph(\n  \"ADD GAME MUSIC\",\n  Word.InsertLocation.before\n);\n\n// 2) New paragraph with \"MAKE  DUMMY MAIN MENU for the Accessibility Settings!\"\n//    where \"MAIN MENU \" is bold.\nconst dummyPara = originalFirstPara.insertParagraph(\n  \"MAKE  DUMMY MAIN MENU for the Accessibility Settings!\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\nconst boldSearch = dummyPara.search(\"MAIN MENU \", { matchCase: true });\nboldSearch.load(\"items\");\nawait context.sync();\nboldSearch.items[0].font.bold = true;\nawait context.sync();\n\n// 3) Rebuild the last paragraph (\"Recognized Voice Commands...\") with the\n//    new wording / extra runs.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst paras = body.paragraphs.items;\nconst lastPara = paras[paras.length - 1];\n\nlastPara.clear();\nawait context.sync();\n\nfunction addRun(text, italic) {\n  const r = lastPara.insertText(text, Word.InsertLocation.end);\n  r.font.name = \"Arial\";\n  r.font.size = 12;\n  if (italic) {\n    r.font.italic = true;\n  }\n  return r;\n}\n\naddRun(\"\\u201cRecognized Voice Commands:\", false);\naddRun(\" YOU SAID\", false);\naddRun(\" \", false);\naddRun(\"\\u2018\", false);\naddRun(\"whatever the user said\", true);\naddRun(\"\\u2019 \", true);\naddRun(\"\\u201d\", false);\n\nawait context.sync();\n", "ps1": "# The document currently starts with \"OBSTACLES: Lerp horizontal and\n# vertical, lerp scale too\" followed by an empty paragraph, etc.\n# This script:\n#   1) Inserts a brand-new paragraph \"ADD GAME MUSIC\" before everything else.\n#   2) Inserts a brand-new paragraph right after it (before the old first\n#      paragraph) containing \"MAKE  DUMMY MAIN MENU for the Accessibility\n#      Settings!\" with \"MAIN MENU \" in bold.\n#   3) Leaves the rest of the body alone, except the final \"Recognized Voice\n#      Commands\" paragraph, whose wording gets extra text spliced in.\n\n$d = $word.ActiveDocument\n\n# --- 1) New paragraph: \"ADD GAME MUSIC\" ------------------------------------\n$firstParaRange = $d.Paragraphs(1).Range\n$firstParaRange.InsertParagraphBefore()\n$titleRange = $d.Paragraphs(1).Range\n$titleRange.Text = \"ADD GAME MUSIC\"\n\n# --- 2) New paragraph: \"MAKE  DUMMY MAIN MENU for the Accessibility Settings!\"\n$oldFirstRange = $d.Paragraphs(2).Range\n$oldFirstRange.InsertParagraphBefore()\n$dummyRange = $d.Paragraphs(2).Range\n$dummyRange.Text = \"MAKE  DUMMY MAIN MENU for the Accessibility Settings!\"\n\n$dummyStart = $d.Paragraphs(2).Range.Start\n$boldStart = $dummyStart + (\"MAKE  DUMMY \").Length\n$boldEnd = $boldStart + (\"MAIN MENU \").Length\n$boldRange = $d.Range($boldStart, $boldEnd)\n$boldRange.Font.Bold = $true\n\n# --- 3) Rebuild the last paragraph (\"Recognized Voice Commands...\") -------\n$leftDouble = [char]0x201C\n$rightDouble = [char]0x201D\n$leftSingle = [char]0x2018\n$rightSingle = [char]0x2019\n\n$count = $d.Paragraphs.Count\n$lastParaRange = $d.Paragraphs($count).Range\n$contentRange = $d.Range($lastParaRange.Start, $lastParaRange.End - 1)\n\n$prefix = $leftDouble + \"Recognized Voice Commands: YOU SAID \" + $leftSingle\n$italicPart = \"whatever the user said\" + $rightSingle + \" \"\n$suffix = $rightDouble\n\n$contentRange.Text = $prefix + $italicPart + $suffix\n\n$pStart = $d.Paragraphs($count).Range.Start\n$italicStart = $pStart + $prefix.Length\n$italicEnd = $italicStart + $italicPart.Length\n$italicRange = $d.Range($italicStart, $italicEnd)\n$italicRange.Font.Italic = $true\n"}
